# Edit script: applies the changes described in the commit diff.
#
# 1. Update all cached "datetimeFigureOut" date-field placeholder text
#    (9/30/2024 -> 10/10/2024) across the slide master and every slide
#    layout (this mirrors PowerPoint automatically refreshing the cached
#    date text for fixed/auto date placeholders whenever the file is
#    re-saved on a later date).
# 2. On slide 2 ("Background & Motivation"):
#    - Re-crop & reposition/resize the picture, and send it to the back
#      of the z-order.
#    - Reposition/resize the existing "Background" content placeholder.
#    - Add a new "Motivation" text box next to it.

$p = $ppt.ActivePresentation

function Update-DateText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -eq -1) {
            $txt = $sh.TextFrame.TextRange.Text
            if ($txt -eq "9/30/2024") {
                $sh.TextFrame.TextRange.Text = "10/10/2024"
            }
        }
    }
}

# --- 1. Refresh cached date-field text everywhere it is reachable ---
$master = $p.SlideMaster
Update-DateText $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DateText $layout.Shapes
}

Update-DateText $p.NotesMaster.Shapes
Update-DateText $p.HandoutMaster.Shapes

# --- 2. Slide 2 content restructuring ---
$slide = $p.Slides.Item(2)

# Shape 1 = title ("标题 1"), Shape 2 = content placeholder ("内容占位符 2"),
# Shape 3 = picture ("图片 7") in the original document order.
$titleSp   = $slide.Shapes.Item(1)
$contentSp = $slide.Shapes.Item(2)
$pic       = $slide.Shapes.Item(3)

# Re-crop the picture.
$pic.PictureFormat.CropLeft  = 130.22676
$pic.PictureFormat.CropTop   = 34.97835
$pic.PictureFormat.CropRight = 71.23176

# Resize & reposition the picture.
$pic.Width  = 193.73125984251968
$pic.Height = 239.5072440944882
$pic.Left   = 766.2687683574803
$pic.Top    = 149.38251968503937

# Send the picture to the back so the text shapes now draw on top of it.
$pic.ZOrder(1)

# Resize & reposition the existing "Background" content placeholder.
$contentSp.Left   = 66.0
$contentSp.Top    = 149.38259842519685
$contentSp.Width  = 368.6088188976378
$contentSp.Height = 336.99236220472443

# Add the new "Motivation" text box.
$newBox = $slide.Shapes.AddTextbox(1, 401.6516571433071, 149.38251968503937, 381.7396087992126, 344.67259842519684)
$newBox.Name = "内容占位符 2"

$tf = $newBox.TextFrame
$tf.MarginLeft   = 7.2
$tf.MarginTop    = 3.6
$tf.MarginRight  = 7.2
$tf.MarginBottom = 3.6

$tr = $tf.TextRange
$tr.Text = "Motivation: "
$tr.Font.Size = 28
[void]$tr.InsertAfter("The motivation mainly stems from the increasingly serious problem of noise pollution in modern society. Noise pollution not only affects people's hearing, sleep and communication, but also accelerates the aging of mechanical structures and buildings. Therefore, the research and development of effective noise reduction materials is of great significance to ensure the sustainable development of economic society.")
$newBox.TextFrame.TextRange.Font.Size = 28

Write-Output "Edit complete"
